$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a literal text value into a cell even when the text
# looks like a plain number (e.g. "568.78"), without leaving any
# number-format/style residue behind -- the cell keeps its original
# style once the value is in place, matching the source workbook where
# every data cell is stored as a literal inline string.
function Set-TextValue($range, $text) {
    $savedStyle = $range.Style
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = $savedStyle
}

$ws.Range('D2').Value = '63.447.29'
$ws.Range('E2').Value = '  +1.44%  '
$ws.Range('D3').Value = '3.408.39'
$ws.Range('E3').Value = '  +2.12%  '
$ws.Range('E4').Value = '  -0.02%  '
Set-TextValue $ws.Range('D5') '568.78'
$ws.Range('E5').Value = '  +1.53%  '
Set-TextValue $ws.Range('D6') '156.45'
$ws.Range('E6').Value = '  +3.29%  '
$ws.Range('E7').Value = '  -0.02%  '
$ws.Range('D8').Value = '3.407.23'
$ws.Range('E8').Value = '  +2.00%  '
Set-TextValue $ws.Range('D9') '0.544'
$ws.Range('E9').Value = '  +2.19%  '
Set-TextValue $ws.Range('D10') '7.38'
$ws.Range('E10').Value = '  -0.17%  '
$ws.Range('E11').Value = '  +3.61%  '
Set-TextValue $ws.Range('D12') '0.430'
$ws.Range('E12').Value = '  -0.85%  '
$ws.Range('D13').Value = '3.995.06'
$ws.Range('E13').Value = '  +2.08%  '
$ws.Range('E14').Value = '  -3.10%  '
$ws.Range('E15').Value = '  +8.15%  '
Set-TextValue $ws.Range('D16') '27.19'
$ws.Range('E16').Value = '  +1.52%  '
$ws.Range('D17').Value = '63.429.22'
$ws.Range('E17').Value = '  +1.41%  '
$ws.Range('D18').Value = '3.416.52'
$ws.Range('E18').Value = '  +1.53%  '
Set-TextValue $ws.Range('D19') '6.23'
$ws.Range('E19').Value = '  -1.65%  '
Set-TextValue $ws.Range('D20') '14.04'
$ws.Range('E20').Value = '  +1.77%  '
Set-TextValue $ws.Range('D21') '377.01'
$ws.Range('E21').Value = '  -1.61%  '
Set-TextValue $ws.Range('D22') '8.06'
$ws.Range('E22').Value = '  -3.92%  '
Set-TextValue $ws.Range('D23') '0.997'
$ws.Range('E23').Value = '  -0.30%  '
Set-TextValue $ws.Range('D24') '71.62'
$ws.Range('E24').Value = '  +2.32%  '
Set-TextValue $ws.Range('D25') '0.527'
$ws.Range('E25').Value = '  -0.87%  '
Set-TextValue $ws.Range('D26') '0.0000120'
$ws.Range('E26').Value = '  +27.00%  '
$ws.Range('E27').Value = '  +4.26%  '
$ws.Range('E28').Value = '  -0.04%  '
Set-TextValue $ws.Range('D29') '0.999'
$ws.Range('E29').Value = '  -0.12%  '
Set-TextValue $ws.Range('D30') '6.05'
$ws.Range('E30').Value = '  +8.46%  '
$ws.Range('B31').Value = 'PancakeSwap'
$ws.Range('C31').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
Set-TextValue $ws.Range('D31') '2.01'
$ws.Range('E31').Value = '  +1.69%  '
$ws.Range('B32').Value = 'Fetch.AI'
$ws.Range('C32').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
Set-TextValue $ws.Range('D32') '1.36'
$ws.Range('E32').Value = '  +4.49%  '
Set-TextValue $ws.Range('D33') '23.21'
$ws.Range('E33').Value = '  +1.56%  '
$ws.Range('B34').Value = 'USDe'
$ws.Range('C34').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
Set-TextValue $ws.Range('D34') '0.998'
$ws.Range('E34').Value = '  +0.00%  '
$ws.Range('B35').Value = 'RenderToken'
$ws.Range('C35').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextValue $ws.Range('D35') '6.33'
$ws.Range('E35').Value = '  -3.30%  '
Set-TextValue $ws.Range('D36') '6.76'
$ws.Range('E36').Value = '  +0.76%  '
Set-TextValue $ws.Range('D37') '158.92'
$ws.Range('E37').Value = '  -0.38%  '
Set-TextValue $ws.Range('D38') '1.44'
$ws.Range('E38').Value = '  -1.66%  '
$ws.Range('D39').Value = '2.976.08'
$ws.Range('E39').Value = '  +6.41%  '
Set-TextValue $ws.Range('D40') '0.0757'
$ws.Range('E40').Value = '  +2.44%  '
Set-TextValue $ws.Range('D41') '27.03'
$ws.Range('E41').Value = '  +0.83%  '
Set-TextValue $ws.Range('D42') '1.81'
$ws.Range('E42').Value = '  -3.67%  '
Set-TextValue $ws.Range('D43') '0.0315'
$ws.Range('E43').Value = '  +0.32%  '
Set-TextValue $ws.Range('D44') '41.89'
$ws.Range('E44').Value = '  +3.63%  '
Set-TextValue $ws.Range('D45') '0.761'
$ws.Range('E45').Value = '  +2.62%  '
Set-TextValue $ws.Range('D46') '4.30'
$ws.Range('E46').Value = '  +1.23%  '
Set-TextValue $ws.Range('D47') '23.19'
$ws.Range('E47').Value = '  +5.85%  '
Set-TextValue $ws.Range('D48') '1.06'
$ws.Range('E48').Value = '  +2.65%  '
Set-TextValue $ws.Range('D49') '2.18'
$ws.Range('E49').Value = '  +22.69%  '
Set-TextValue $ws.Range('D50') '6.32'
$ws.Range('E50').Value = '  +0.61%  '
Set-TextValue $ws.Range('D51') '293.92'
$ws.Range('E51').Value = '  +2.08%  '
